$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "who is the performer associated with green?`n?"
$ws.Range("H2").Value = "JIMMY SAVILE"
$ws.Range("R2").Value = 0.3076923191547394

$ws.Range("E3").Value = "who is the spouse of steve hillage?`n?"
$ws.Range("H3").Value = "MARY ANNE HOBBS"
$ws.Range("J3").Value = 0
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = 0
$ws.Range("N3").Value = 0.1
